$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Diaria" rows pulled from the MV data source: date label in column A
# (kept as plain text, matching the existing date-label cells above it),
# 10000 in column B (Cupo) and 0 in column D (Total monto adjudicado).
$dates = @("08-09-2021", "09-09-2021", "14-09-2021", "15-09-2021", "16-09-2021")

$row = 20
foreach ($d in $dates) {
    $cell = $ws.Cells.Item($row, 1)

    # Typing a dd-mm-yyyy-looking string straight into a General-formatted
    # cell would get auto-recognised as a date serial. Write it as a text
    # formula first (guaranteed string result) and then flatten it to a
    # plain value via copy / paste-values, which preserves the literal text
    # without touching the cell's number format/style.
    $cell.Formula = '="' + $d + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = 10000
    $ws.Cells.Item($row, 4).Value = 0

    $row++
}
